$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 546
$ws1.Range("F9").Value = 408
$ws1.Range("F10").Value = 3464
$ws1.Range("F11").Value = 49

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 546
$ws4.Range("F10").Value = 408
$ws4.Range("F11").Value = 3464
$ws4.Range("F12").Value = 49
